# Regenerate the localization status report for archive:
# flip every "Ready for handoff" status cell over to "In Translation"
# across the Overview summary sheet and each per-language detail sheet,
# then let the Status column narrow to fit the new (shorter) text, just
# like Excel does when a table column's content shrinks.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn / de-de status columns (E, F) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Cells.Replace("Ready for handoff", "In Translation")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn detail sheet: Status column (C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Cells.Replace("Ready for handoff", "In Translation")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de detail sheet: Status column (C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Cells.Replace("Ready for handoff", "In Translation")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
